$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the last existing data row (33) three times, preserving all
# formatting (number formats, styles, etc.) exactly as Excel would when a
# user copies a row and inserts copies below it.
$ws.Rows(33).Copy()
$ws.Rows(34).Insert()
$ws.Rows(33).Copy()
$ws.Rows(34).Insert()
$ws.Rows(33).Copy()
$ws.Rows(34).Insert()

# New user records added as part of the 16th May refresh.
$newUsers = @(
    @{ id = 110033; uin = 9317596771; name = "Nikola Tesla"; email = "nikola.tesla@xyz.com"; mobile = 818876434 },
    @{ id = 110034; uin = 9317596772; name = "Graham Bell";  email = "graham.bell@xyz.com";  mobile = 818876435 },
    @{ id = 110035; uin = 9317596773; name = "Albert Miles"; email = "albert.miles@xyz.com"; mobile = 818876436 }
)

# Populate column-by-column (all ids, then all uins, then all names, ...)
# so that new shared-string entries are interned in the same order Excel
# produced them in (names first, then emails).
for ($i = 0; $i -lt $newUsers.Count; $i++) {
    $ws.Cells.Item(34 + $i, 1).Value = $newUsers[$i].id
}
for ($i = 0; $i -lt $newUsers.Count; $i++) {
    $ws.Cells.Item(34 + $i, 2).Value = $newUsers[$i].uin
}
for ($i = 0; $i -lt $newUsers.Count; $i++) {
    $ws.Cells.Item(34 + $i, 3).Value = $newUsers[$i].name
}
for ($i = 0; $i -lt $newUsers.Count; $i++) {
    $ws.Cells.Item(34 + $i, 4).Value = $newUsers[$i].email
}
for ($i = 0; $i -lt $newUsers.Count; $i++) {
    $ws.Cells.Item(34 + $i, 5).Value = $newUsers[$i].mobile
}

# Reflect the new used range / active selection like Excel would after
# adding rows and then selecting the entire row band below the data
# (mirrors the original "select rest of sheet below data" state).
$ws.Range("A37:XFD1048576").Select()
